$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "750×6="; New = "976×7=" },
    @{ Old = "844×5="; New = "744×9=" },
    @{ Old = "464×3="; New = "830×6=" },
    @{ Old = "374×6="; New = "777×3=" },
    @{ Old = "668×2="; New = "521×7=" },
    @{ Old = "897×2="; New = "240×6=" },
    @{ Old = "306×3="; New = "722×4=" },
    @{ Old = "170×5="; New = "986×5=" },
    @{ Old = "754×6="; New = "716×4=" },
    @{ Old = "831×3="; New = "294×7=" },
    @{ Old = "565×4="; New = "933×4=" },
    @{ Old = "273×4="; New = "404×9=" },
    @{ Old = "492×2="; New = "308×7=" },
    @{ Old = "582×4="; New = "930×9=" },
    @{ Old = "742×8="; New = "953×4=" },
    @{ Old = "819×7="; New = "246×3=" },
    @{ Old = "354×3="; New = "937×4=" },
    @{ Old = "567×5="; New = "941×2=" },
    @{ Old = "920×3="; New = "965×3=" },
    @{ Old = "117×5="; New = "405×4=" },
    @{ Old = "173×2="; New = "210×5=" },
    @{ Old = "224×9="; New = "441×4=" },
    @{ Old = "128×9="; New = "136×3=" },
    @{ Old = "173×5="; New = "441×2=" },
    @{ Old = "474×2="; New = "710×5=" }
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.Old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $r.New, 2)
}
